# Remove column M (the old "days since last drink" helper column) from the
# alcohol measurement sheet, shifting column N left to take its place.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("M").Delete()
$ws.Range("M1").Select()
